# Remove TOPINS for HETD (included directly in RSD)
#
# The ~TFM_TOPINS table on the PWR sheet had a single data row (row 4)
# mapping FT-RSDHET -> HETD via an "IN" insertion. Since HETD is now
# included directly in RSD, that row is no longer needed, so delete it
# entirely. Everything below it (the ~TFM_INS table and its data rows)
# shifts up by one row automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PWR")

# Delete the whole row 4 (F4:I4 = IN / IN / FT-RSDHET / HETD)
$ws.Rows("4:4").Delete()

# Make PWR the active sheet/tab and move the selection to J17,
# matching where the editor ended up after the edit.
$ws.Activate()
$ws.Range("J17").Select()
